$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1006620.5
$ws.Range("I62").Value = 2504098.8
$ws.Range("J62").Value = 8301.666999999999
$ws.Range("K62").Value = 2504098.8
$ws.Range("L62").Value = 8301.666999999999
$ws.Range("M62").Value = -2503474.8
$ws.Range("N62").Value = -9549.666999999999

$ws.Range("H65").Value = 1006620.5
$ws.Range("I65").Value = 2504098.8
$ws.Range("J65").Value = 8301.666999999999
$ws.Range("K65").Value = 12520494
$ws.Range("L65").Value = 41508.335
$ws.Range("M65").Value = -12517374
$ws.Range("N65").Value = -47748.335

$ws.Range("H82").Value = 2360.5
$ws.Range("I82").Value = 2360.5
$ws.Range("K82").Value = 7081.5
$ws.Range("M82").Value = -6675.5

$ws.Range("H85").Value = 2360.5
$ws.Range("I85").Value = 2360.5
$ws.Range("K85").Value = 7081.5
$ws.Range("M85").Value = -5677.5

$ws.Range("H106").Value = 15875473
$ws.Range("I106").Value = 15875473
$ws.Range("K106").Value = 15875473
$ws.Range("M106").Value = -15874842

$ws.Range("H129").Value = 1001.04083
$ws.Range("J129").Value = 1057.5778
$ws.Range("L129").Value = 3172.7334
$ws.Range("N129").Value = -13172.7334

$ws.Range("H138").Value = 1719.18
$ws.Range("I138").Value = 778.5
$ws.Range("J138").Value = 2049.6892
$ws.Range("K138").Value = 2335.5
$ws.Range("L138").Value = 6149.067599999999
$ws.Range("M138").Value = 2804.5
$ws.Range("N138").Value = -16429.0676


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 200069
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 200069
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 200069
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -205561


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 262.25
$ws.Range("I22").Value = 262.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 262.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -89.25
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 1008.1053
$ws.Range("I94").Value = 1089.0588
$ws.Range("K94").Value = 1089.0588
$ws.Range("M94").Value = -638.0588

$ws.Range("H134").Value = 3696.182
$ws.Range("I134").Value = 2783.8572
$ws.Range("J134").Value = 5292.75
$ws.Range("K134").Value = 8351.571599999999
$ws.Range("L134").Value = 15878.25
$ws.Range("M134").Value = -5816.571599999999
$ws.Range("N134").Value = -20948.25


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7642.6665
$ws.Range("I31").Value = 3000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2705

$ws.Range("H34").Value = 7642.6665
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2798


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 6000
$ws.Range("J17").Value = 6000
$ws.Range("L17").Value = 18000
$ws.Range("N17").Value = -18338

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H41").Value = 650
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 800
$ws.Range("K41").Value = 600
$ws.Range("L41").Value = 2400
$ws.Range("M41").Value = -262
$ws.Range("N41").Value = -3076

$ws.Range("H58").Value = 2250
$ws.Range("I58").Value = 2333.3333
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 6999.999899999999
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -6871.999899999999
$ws.Range("N58").Value = -6256

$ws.Range("H64").Value = 2996.6667
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 3995
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 11985
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -12525

$ws.Range("H67").Value = 2996.6667
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 3995
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 11985
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -13857

$ws.Range("H68").Value = 15631.25
$ws.Range("I68").Value = 24540.4
$ws.Range("J68").Value = 782.6667
$ws.Range("K68").Value = 73621.20000000001
$ws.Range("L68").Value = 2348.0001
$ws.Range("M68").Value = -72810.20000000001
$ws.Range("N68").Value = -3970.0001

$ws.Range("H70").Value = 2324.7778
$ws.Range("I70").Value = 784.6
$ws.Range("K70").Value = 2353.8
$ws.Range("M70").Value = -2038.8

$ws.Range("H71").Value = 15631.25
$ws.Range("I71").Value = 24540.4
$ws.Range("J71").Value = 782.6667
$ws.Range("K71").Value = 220863.6
$ws.Range("L71").Value = 7044.0003
$ws.Range("M71").Value = -216807.6
$ws.Range("N71").Value = -15156.0003

$ws.Range("H73").Value = 2324.7778
$ws.Range("I73").Value = 784.6
$ws.Range("K73").Value = 2353.8
$ws.Range("M73").Value = -1261.8

$ws.Range("H94").Value = 2076.2307
$ws.Range("I94").Value = 918.2
$ws.Range("J94").Value = 2800
$ws.Range("K94").Value = 2754.6
$ws.Range("L94").Value = 8400
$ws.Range("M94").Value = -2078.6
$ws.Range("N94").Value = -9752

$ws.Range("H106").Value = 3400
$ws.Range("J106").Value = 3400
$ws.Range("L106").Value = 10200
$ws.Range("N106").Value = -12092

$ws.Range("H109").Value = 2666.6667
$ws.Range("I109").Value = 1500
$ws.Range("J109").Value = 3250
$ws.Range("K109").Value = 4500
$ws.Range("L109").Value = 9750
$ws.Range("M109").Value = -3460
$ws.Range("N109").Value = -11830

$ws.Range("H112").Value = 333334270
$ws.Range("I112").Value = 1413.5
$ws.Range("K112").Value = 4240.5
$ws.Range("M112").Value = -3132.5

$ws.Range("H124").Value = 1033.3334
$ws.Range("I124").Value = 750
$ws.Range("J124").Value = 1600
$ws.Range("K124").Value = 2250
$ws.Range("L124").Value = 4800
$ws.Range("M124").Value = 2660
$ws.Range("N124").Value = -14620

$ws.Range("H129").Value = 1131.8334
$ws.Range("I129").Value = 562.5
$ws.Range("J129").Value = 1416.5
$ws.Range("K129").Value = 1687.5
$ws.Range("L129").Value = 4249.5
$ws.Range("M129").Value = 3312.5
$ws.Range("N129").Value = -14249.5

$ws.Range("H131").Value = 2384.489
$ws.Range("J131").Value = 2531.4302
$ws.Range("L131").Value = 7594.290599999999
$ws.Range("N131").Value = -17674.2906


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7305.5557
$ws.Range("I70").Value = 8353.846
$ws.Range("J70").Value = 4580
$ws.Range("K70").Value = 8353.846
$ws.Range("L70").Value = 4580
$ws.Range("M70").Value = -8083.846
$ws.Range("N70").Value = -5120

$ws.Range("H73").Value = 7305.5557
$ws.Range("I73").Value = 8353.846
$ws.Range("J73").Value = 4580
$ws.Range("K73").Value = 8353.846
$ws.Range("L73").Value = 4580
$ws.Range("M73").Value = -7417.846
$ws.Range("N73").Value = -6452

$ws.Range("H80").Value = 2837.5
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -5246

$ws.Range("H83").Value = 2837.5
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -26234

$ws.Range("H126").Value = 2190.1538
$ws.Range("I126").Value = 1869.1333
$ws.Range("J126").Value = 2390.7917
$ws.Range("K126").Value = 5607.3999
$ws.Range("L126").Value = 7172.375100000001
$ws.Range("M126").Value = -3137.3999
$ws.Range("N126").Value = -12112.3751

$ws.Range("H132").Value = 3187.6
$ws.Range("I132").Value = 2870.6206
$ws.Range("K132").Value = 8611.861800000001
$ws.Range("M132").Value = -6081.861800000001

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16669170
$ws.Range("I132").Value = 20835402
$ws.Range("K132").Value = 62506206
$ws.Range("M132").Value = -62503676

